$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new item row at row 14 (pushes old row14/15 down to 15/16) ---
$ws.Range("A14:Q14").Insert(-4121)

# Copy the formatting from the row above (row 13, an existing item row) so the
# new row matches the look of the other item rows (border/fill/font/numfmt).
$ws.Range("A13:Q13").Copy()
$ws.Range("A14:Q14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the merged cells for the new row, matching the other item rows.
$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

# Row height matches the alternating pattern used by the item rows.
$ws.Rows(14).RowHeight = 25.5

# --- Fill in the values for the new item (#8) ---
$ws.Range("A14").Value2 = 8
$ws.Range("C14").Value2 = "صابون ديتول العنايه بالبشره"
$ws.Range("H14").Value2 = "13:0"

# L14 and P14 keep number-style formatting but must store their numeric-looking
# value as text (like the other rows do) - force text entry, then restore the
# original number format so the cell's visual style is unchanged.
$fmtL = $ws.Range("L14").NumberFormat
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value2 = "0"
$ws.Range("L14").NumberFormat = $fmtL

$ws.Range("N14").Value2 = "45.00"

$fmtP = $ws.Range("P14").NumberFormat
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value2 = "45.0000"
$ws.Range("P14").NumberFormat = $fmtP

$ws.Range("Q14").Value2 = "1:0"

# --- Update the totals row (now shifted to row 15) ---
$ws.Rows(15).RowHeight = 24.75
$ws.Range("P15").Value2 = 572
